$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.921.55"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -0.41%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'1.907.98"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -0.19%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'0.9992"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.66%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'313.43"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -0.64%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'0.9984"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -0.60%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.4992"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  +3.68%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.3815"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  +0.14%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.07298"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -0.86%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'0.9111"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -2.46%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'21.11"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +1.30%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'0.07690"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -1.16%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'1.889.07"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -2.28%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'5.507"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +0.14%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'92.66"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +0.67%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'0.9994"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -0.66%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'0.000008737"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -1.48%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'0.9989"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -0.53%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'27.944.12"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  -0.45%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'14.65"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -0.85%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'5.182"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +0.32%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = "'  -0.64%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'6.584"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -1.01%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').Value = "'  -1.81%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'1.856"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -3.44%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'2.221"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +4.25%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = "'  -0.49%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'115.42"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -1.30%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'4.902"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -1.31%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'0.09023"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  +0.74%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'3.202"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -3.12%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'4.873"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +4.20%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'  -2.65%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'0.7722"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -0.95%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'0.02090"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  +1.62%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'2.559"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -1.96%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'3.068"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  +2.31%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'  -1.69%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'0.5554"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E40').Value = "'  -0.70%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'6.882"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -2.13%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'8.506"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +0.23%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('B43').Value = "'Quant"
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').Value = "'112.74"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +3.88%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('B44').Value = "'Algorand"
$ws.Range('B44').Style = 'Normal'
$ws.Range('C44').Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range('C44').Style = 'Normal'
$ws.Range('D44').Value = "'0.1521"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -0.46%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Value = "'  -0.37%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'0.4839"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +0.21%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'0.9983"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -0.60%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'1.637"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -0.72%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'67.43"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -1.14%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'  -0.51%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'0.9072"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +0.77%  "
$ws.Range('E51').Style = 'Normal'
